$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The target content for rows 2-4 duplicates content that already exists
# further down the sheet (row 3 -> row 2, row 9 -> row 3, row 10 -> row 4).
# Use Range.Copy(destination) instead of re-typing the values: column B
# holds numeric-looking text (e.g. "8020023215") that must stay text, and
# copying the cell (rather than assigning a freshly-typed string to .Value)
# preserves the original text type without Excel's "looks like a number"
# auto-conversion and without introducing a new cell style.

# Row 2 gets row 3's current (pre-edit) content - read it out first since
# row 3 itself is about to be overwritten below.
$ws.Range("A3:C3").Copy($ws.Range("A2:C2"))

# Row 3 gets row 9's content.
$ws.Range("A9:C9").Copy($ws.Range("A3:C3"))

# Row 4 gets row 10's content.
$ws.Range("A10:C10").Copy($ws.Range("A4:C4"))
